# Commit message: "Fruta / hortaliza, semanal" (weekly fruit/vegetable price update)
#
# The diff shows a brand-new price-report row inserted above the existing
# row 305 (Hortaliza / Femacal de La Calera / Poroto verde sheet), pushing
# every following row down by one (old row 305 -> new row 306, ...,
# old row 382 -> new row 383), and the sheet's dimension growing from
# A1:R382 to A1:R383.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 305; Excel shifts rows 305:382 down to 306:383.
$ws.Rows.Item(305).Insert()

# Populate the newly inserted row 305 with the new weekly observation.
$ws.Cells.Item(305, 1).Value  = 3
$ws.Cells.Item(305, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(305, 3).Value  = "Coquimbo"
$ws.Cells.Item(305, 4).Value  = 44722
$ws.Cells.Item(305, 5).Value  = 5
$ws.Cells.Item(305, 6).Value  = 100112031
$ws.Cells.Item(305, 7).Value  = "Poroto verde"
$ws.Cells.Item(305, 8).Value  = "Magnum"
$ws.Cells.Item(305, 9).Value  = "Primera"
$ws.Cells.Item(305, 10).Value = 85
$ws.Cells.Item(305, 11).Value = 28000
$ws.Cells.Item(305, 12).Value = 29000
$ws.Cells.Item(305, 13).Value = 28471
$ws.Cells.Item(305, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(305, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(305, 16).Value = 1139
$ws.Cells.Item(305, 17).Value = 25
$ws.Cells.Item(305, 18).Value = "Hortaliza"
